$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: new survey response (Isabel Urdiales Sotres) ---

# Do the "forced text" cells FIRST (SmartScore "0.000" text cells and the blank
# B15 cell): they rely on a temporary text number-format + ClearFormats() to land
# on the default (unstyled) cell style. Doing this before the multi-line G15 write
# keeps a later ClearFormats() from picking up an incidental "best fit" font that
# Excel associates with the row once it contains wrapped/multi-line text.

$smartScoreTextCells = @("I15", "L15", "O15", "R15", "U15", "X15", "AA15", "AD15", "AG15")
foreach ($cellRef in $smartScoreTextCells) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = "0.000"
    $cell.ClearFormats()
}

# B15 (Grupo_Experimental) is blank for this response - create the cell explicitly
# (otherwise no cell is emitted at all) while keeping it empty and unstyled.
$bCell = $ws.Range("B15")
$bCell.NumberFormat = "@"
$bCell.Value = ""
$bCell.ClearFormats()

# Remaining plain text / numeric cells

$ws.Range("D15").Value = 20
$ws.Range("A15").Value = @"
Isabel Urdiales Sotres_20251120_223408
"@
$ws.Range("C15").Value = @"
Isabel Urdiales Sotres
"@
$ws.Range("E15").Value = @"
Female
"@
$ws.Range("F15").Value = @"
2025-11-20 22:34:08
"@
$ws.Range("G15").Value = @"
{
  "portion": 0.0,
  "diet": 0.0,
  "salt": 0.0,
  "fat": 0.0,
  "natural": 0.0,
  "convenience": 0.0,
  "price": 0.0
}
"@
$ws.Range("H15").Value = @"
Maruchan Ramen Sabor Pollo
"@
$ws.Range("J15").Value = @"
Sabor clásico, económico, alto en sodio, no saludable, nostálgico
"@
$ws.Range("K15").Value = @"
Nissin Chow Mein Teriyaki Beef
"@
$ws.Range("M15").Value = @"
Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa
"@
$ws.Range("N15").Value = @"
Nongshim Shin Ramyun
"@
$ws.Range("P15").Value = @"
Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio
"@
$ws.Range("Q15").Value = @"
Annie’s Shells & White Cheddar
"@
$ws.Range("S15").Value = @"
Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños
"@
$ws.Range("T15").Value = @"
Velveeta Original Shells & Cheese (microwave cups)
"@
$ws.Range("V15").Value = @"
Muy cremoso, porción individual, rápido, salado, ideal para niños
"@
$ws.Range("W15").Value = @"
Kraft Macaroni & Cheese Dinner
"@
$ws.Range("Y15").Value = @"
Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato
"@
$ws.Range("Z15").Value = @"
Wild Planet Wild Tuna Pasta Salad
"@
$ws.Range("AB15").Value = @"
Sabor fresco, buena proteína, saludable, porción algo pequeña
"@
$ws.Range("AC15").Value = @"
StarKist Chicken Creations (Chicken Salad)
"@
$ws.Range("AE15").Value = @"
Portátil, saludable, fácil, buena textura, sabor suave
"@
$ws.Range("AF15").Value = @"
Kitchens of India Variety Pack
"@
$ws.Range("AH15").Value = @"
Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad
"@

# The multi-line "Pesos" JSON text in G15 makes Excel auto-expand the row height;
# auto-fit it back down so row 15 matches the other (un-customized) rows.
$ws.Rows.Item(15).EntireRow.AutoFit()
